# Updated cryptos list on Sat Jun  3 13:31:38 UTC 2023 with GitHub Actions
#
# The workbook's Coin/Link/Price/Volume(1h) columns (B:E) hold text values
# even when they look numeric (European-style "27.200.50" price strings,
# "1.002" style prices, "  +0.39%  " volume strings, etc). Excel's COM
# `.Value` setter auto-coerces anything that parses as a number, which would
# both change the stored type (t="inlineStr"/shared-string -> numeric <v>)
# and -- worse -- mangle exact text like trailing zeros ("0.08460" ->
# 0.0846). To keep every touched cell byte-for-byte text, force the cell to
# Text format, write via Value2, then drop the format back to the sheet's
# normal/default style so we don't leave a stray number-format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "27.205.30"

# Row 3 - Ethereum
Set-TextValue $ws "D3" "1.905.49"
Set-TextValue $ws "E3" "  +0.77%  "

# Row 4 - TetherUSD
Set-TextValue $ws "D4" "1.002"

# Row 5 - BNB
Set-TextValue $ws "D5" "307.86"
Set-TextValue $ws "E5" "  +0.39%  "

# Row 6 - USDC
Set-TextValue $ws "D6" "1.002"
Set-TextValue $ws "E6" "  +0.12%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.5196"
Set-TextValue $ws "E7" "  +0.77%  "

# Row 8 - Cardano
Set-TextValue $ws "D8" "0.3765"
Set-TextValue $ws "E8" "  +0.14%  "

# Row 9 - Dogecoin
Set-TextValue $ws "D9" "0.07269"
Set-TextValue $ws "E9" "  +1.02%  "

# Row 10 - Solana
Set-TextValue $ws "D10" "21.18"
Set-TextValue $ws "E10" "  +0.21%  "

# Row 11 - Polygon
Set-TextValue $ws "E11" "  +0.29%  "

# Row 12 - TRON
Set-TextValue $ws "D12" "0.08460"

# Row 13 - was Litecoin, now WrappedEther
Set-TextValue $ws "B13" "WrappedEther"
Set-TextValue $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D13" "1.919.03"
Set-TextValue $ws "E13" "  +1.51%  "

# Row 14 - was WrappedEther, now Litecoin
Set-TextValue $ws "B14" "Litecoin"
Set-TextValue $ws "C14" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D14" "96.78"
Set-TextValue $ws "E14" "  +2.59%  "

# Row 15 - Polkadot
Set-TextValue $ws "D15" "5.297"
Set-TextValue $ws "E15" "  +0.99%  "

# Row 16 - BinanceUSD
Set-TextValue $ws "D16" "1.002"
Set-TextValue $ws "E16" "  +0.13%  "

# Row 17 - ShibaInu
Set-TextValue $ws "D17" "0.000008657"
Set-TextValue $ws "E17" "  +1.88%  "

# Row 18 - Avalanche
Set-TextValue $ws "E18" "  +0.81%  "

# Row 19 - Dai
Set-TextValue $ws "E19" "  +0.13%  "

# Row 20 - WrappedBTC
Set-TextValue $ws "D20" "27.245.82"
Set-TextValue $ws "E20" "  +0.42%  "

# Row 21 - Uniswap
Set-TextValue $ws "D21" "5.096"
Set-TextValue $ws "E21" "  +0.72%  "

# Row 22 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D22" "2.147.99"
Set-TextValue $ws "E22" "  +0.43%  "

# Row 23 - Cosmos
Set-TextValue $ws "D23" "10.66"
Set-TextValue $ws "E23" "  +0.81%  "

# Row 24 - Chainlink
Set-TextValue $ws "D24" "6.441"

# Row 25 - LidoDAOToken
Set-TextValue $ws "D25" "2.345"
Set-TextValue $ws "E25" "  +2.71%  "

# Row 26 - Monero
Set-TextValue $ws "D26" "147.06"
Set-TextValue $ws "E26" "  +0.31%  "

# Row 27 - Toncoin
Set-TextValue $ws "D27" "1.754"
Set-TextValue $ws "E27" "  +0.84%  "

# Row 28 - EthereumClassic
Set-TextValue $ws "D28" "18.24"
Set-TextValue $ws "E28" "  +0.94%  "

# Row 29 - BitcoinCash
Set-TextValue $ws "D29" "115.15"
Set-TextValue $ws "E29" "  +0.58%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue $ws "E30" "  +0.45%  "

# Row 31 - Filecoin
Set-TextValue $ws "D31" "4.906"
Set-TextValue $ws "E31" "  -0.48%  "

# Row 32 - Stellar
Set-TextValue $ws "D32" "0.09283"
Set-TextValue $ws "E32" "  +0.88%  "

# Row 33 - Hedera
Set-TextValue $ws "D33" "0.05069"
Set-TextValue $ws "E33" "  +0.46%  "

# Row 34 - ImmutableX
Set-TextValue $ws "D34" "0.7958"
Set-TextValue $ws "E34" "  +3.63%  "

# Row 35 - ARBITRUM
Set-TextValue $ws "D35" "1.243"
Set-TextValue $ws "E35" "  +0.59%  "

# Row 36 - MXToken
Set-TextValue $ws "D36" "3.433"
Set-TextValue $ws "E36" "  +4.55%  "

# Row 37 - HuobiToken
Set-TextValue $ws "D37" "2.951"
Set-TextValue $ws "E37" "  -1.38%  "

# Row 38 - TheSandbox
Set-TextValue $ws "D38" "0.5822"
Set-TextValue $ws "E38" "  +3.85%  "

# Row 39 - RenderToken
Set-TextValue $ws "D39" "2.584"
Set-TextValue $ws "E39" "  -0.39%  "

# Row 40 - VeChain
Set-TextValue $ws "D40" "0.02005"
Set-TextValue $ws "E40" "  +0.72%  "

# Row 41 - TrustWalletToken
Set-TextValue $ws "E41" "  +0.36%  "

# Row 42 - Aptos
Set-TextValue $ws "D42" "9.074"
Set-TextValue $ws "E42" "  +0.03%  "

# Row 43 - FraxShare
Set-TextValue $ws "D43" "6.609"
Set-TextValue $ws "E43" "  -0.47%  "

# Row 44 - Quant
Set-TextValue $ws "D44" "116.56"
Set-TextValue $ws "E44" "  -1.67%  "

# Row 45 - Algorand
Set-TextValue $ws "D45" "0.1520"
Set-TextValue $ws "E45" "  +1.17%  "

# Row 46 - Decentraland
Set-TextValue $ws "D46" "0.4885"
Set-TextValue $ws "E46" "  +1.36%  "

# Row 47 - was EnergySwap, now PaxDollar
Set-TextValue $ws "B47" "PaxDollar"
Set-TextValue $ws "C47" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws "D47" "1.002"
Set-TextValue $ws "E47" "  +0.13%  "

# Row 48 - was PaxDollar, now EnergySwap
Set-TextValue $ws "B48" "EnergySwap"
Set-TextValue $ws "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D48" "10.16"
Set-TextValue $ws "E48" "  +0.34%  "

# Row 49 - NEARProtocol
Set-TextValue $ws "D49" "1.636"
Set-TextValue $ws "E49" "  +2.43%  "

# Row 50 - Elrond
Set-TextValue $ws "D50" "37.65"
Set-TextValue $ws "E50" "  +0.08%  "

# Row 51 - Aave
Set-TextValue $ws "D51" "64.05"
Set-TextValue $ws "E51" "  +0.08%  "
